# Revised Budget.xlsx update — 6/19/2017 additions
# (radio adapter / ethernet cable line items), IP Camera item rename,
# and a refreshed Total row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# 1. Make room: insert 4 blank rows right after the old "Total" row (8),
#    pushing the old trailing filler rows (9-24) down to (13-28).
# ---------------------------------------------------------------------
$ws.Rows("9:12").Insert()

# ---------------------------------------------------------------------
# 2. Rename the IP Camera line item in row 7.
# ---------------------------------------------------------------------
$ws.Range("A7").Value = "IP Camera "

# ---------------------------------------------------------------------
# 3. The old row 8 was the "Total" row (B8="Total", E8=SUM(E2:E7)).
#    It becomes the new "ADDITIONS" section header for the 6/19/2017 batch.
# ---------------------------------------------------------------------
$ws.Range("B8").ClearContents()
$ws.Range("E8").ClearContents()
$ws.Range("A8").Value = "ADDITIONS: made on 6/19/2017"
$ws.Range("A8").Font.Bold = $true

# ---------------------------------------------------------------------
# 4. New line item - row 9: 6" Ethernet cables.
# ---------------------------------------------------------------------
$ws.Range("A9").Value = "6 Inch Ethernet Cables"
$ws.Range("B9").Value = "5-Pack 6-inch CAT6 Network UTP Ethernet RJ45"
$ws.Range("C9").Value = 1
$ws.Range("D9").Value = 11.79
$ws.Range("E9").Formula = "=C9*D9"
$ws.Range("F9").Value = "https://www.amazon.com/CablesOnline-Network-Ethernet-Flat-Design-U6-000FK-5/dp/B00OJYZMFG/ref=sr_1_5?ie=UTF8&qid=1497639766&sr=8-5&keywords=6+inch+ethernet+cable"

# ---------------------------------------------------------------------
# 5. New line item - row 10: radio adapter.
# ---------------------------------------------------------------------
$ws.Range("A10").Value = "Radio Adapter"
$ws.Range("B10").Value = "N-Male to SMA Female Adapter"
$ws.Range("C10").Value = 3
$ws.Range("D10").Value = 4.99
$ws.Range("E10").Formula = "=C10*D10"
$ws.Range("F10").Value = "https://www.readymaderc.com/store/index.php?main_page=product_info&cPath=11_45_58&products_id=473"

# Row 11 stays blank (spacer row between the additions and the new Total).

# ---------------------------------------------------------------------
# 6. New "Total" row, now at row 12, summing everything above it.
# ---------------------------------------------------------------------
$ws.Range("B12").Value = "Total"
$ws.Range("E12").Formula = "=SUM(E2:E11)"

# ---------------------------------------------------------------------
# 7. Leave the selection where the editor last left off.
# ---------------------------------------------------------------------
$ws.Range("C11").Select() | Out-Null

Write-Output "Applied 6/19/2017 budget additions."
